$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78, shifting existing rows 78-159 down to 79-160.
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new record's data.
$ws.Cells.Item(78, 1).Value = 5
$ws.Cells.Item(78, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(78, 3).Value = "Maule"
$ws.Cells.Item(78, 4).Value = 44546
$ws.Cells.Item(78, 5).Value = 7
$ws.Cells.Item(78, 6).Value = 100112030
$ws.Cells.Item(78, 7).Value = "Poroto granado"
$ws.Cells.Item(78, 8).Value = "Sin especificar"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 50
$ws.Cells.Item(78, 11).Value = 42000
$ws.Cells.Item(78, 12).Value = 42000
$ws.Cells.Item(78, 13).Value = 42000
$ws.Cells.Item(78, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(78, 15).Value = "Región del Maule"
$ws.Cells.Item(78, 16).Value = 1680
$ws.Cells.Item(78, 17).Value = 25
$ws.Cells.Item(78, 18).Value = "Hortaliza"
